$d = $word.ActiveDocument

# Split "The end." into three runs: "The", " real", " end." by inserting
# " real" right after "The" and toggling a formatting property on/off so
# the engine doesn't re-merge the newly inserted run with its neighbours.
$findRng = $d.Content
$findRng.Find.Text = "The"
$findRng.Find.Forward = $true
$findRng.Find.Execute() | Out-Null
$findRng.Collapse(0)
$findRng.InsertAfter(" real")
$findRng.Bold = 1
$findRng.Bold = 0

# Move the "_GoBack" bookmark from the "...or is it?" paragraph down to the
# very last (empty) paragraph of the document. Word keeps a single
# "_GoBack" bookmark, so adding a new one with the same name removes the
# old one automatically.
$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
